$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the condition text in F7 from "autogen == $param" to "autogen = $param"
$ws.Range("F7").Value = "autogen = `$param"

# Match the saved selection/active cell state (F7) recorded in the sheet view
$ws.Range("F7").Select()
